$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "~1"
$ws.Range("B2").Value = "~1"
$ws.Range("C2").Value = "~1"
$ws.Range("D2").Value = "~1"
$ws.Range("E2").Value = "Phi(~1)p(~1)pent(~1)N(~1)"
$ws.Range("F2").Value = 4
$ws.Range("G2").Value = 87.04659848148148
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0.3330270270279794
$ws.Range("J2").Value = -41.443963

$ws.Range("A3").Value = "~1"
$ws.Range("B3").Value = "~temp"
$ws.Range("C3").Value = "~1"
$ws.Range("D3").Value = "~1"
$ws.Range("E3").Value = "Phi(~1)p(~temp)pent(~1)N(~1)"
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = 87.1192323076923
$ws.Range("H3").Value = 0.0726338262108186
$ws.Range("I3").Value = 0.3211494964769983
$ws.Range("J3").Value = -44.19754

$ws.Range("A4").Value = "~1"
$ws.Range("B4").Value = "~sundur"
$ws.Range("C4").Value = "~1"
$ws.Range("D4").Value = "~1"
$ws.Range("E4").Value = "Phi(~1)p(~sundur)pent(~1)N(~1)"
$ws.Range("F4").Value = 5
$ws.Range("G4").Value = 89.86769930769231
$ws.Range("H4").Value = 2.821100826210824
$ws.Range("I4").Value = 0.08126157201169604
$ws.Range("J4").Value = -41.449073

$ws.Range("A5").Value = "~1"
$ws.Range("B5").Value = "~wind"
$ws.Range("C5").Value = "~1"
$ws.Range("D5").Value = "~1"
$ws.Range("E5").Value = "Phi(~1)p(~wind)pent(~1)N(~1)"
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 89.87052730769231
$ws.Range("H5").Value = 2.823928826210832
$ws.Range("I5").Value = 0.081146749347626
$ws.Range("J5").Value = -41.446245

$ws.Range("A6").Value = "~1"
$ws.Range("B6").Value = "~temp + wind"
$ws.Range("C6").Value = "~1"
$ws.Range("D6").Value = "~1"
$ws.Range("E6").Value = "Phi(~1)p(~temp + wind)pent(~1)N(~1)"
$ws.Range("F6").Value = 6
$ws.Range("G6").Value = 89.94076099999999
$ws.Range("H6").Value = 2.894162518518513
$ws.Range("I6").Value = 0.07834658564970837
$ws.Range("J6").Value = -44.428319

$ws.Range("A7").Value = "~1"
$ws.Range("B7").Value = "~temp + sundur"
$ws.Range("C7").Value = "~1"
$ws.Range("D7").Value = "~1"
$ws.Range("E7").Value = "Phi(~1)p(~temp + sundur)pent(~1)N(~1)"
$ws.Range("F7").Value = 6
$ws.Range("G7").Value = 90.13548400000001
$ws.Range("H7").Value = 3.088885518518524
$ws.Range("I7").Value = 0.07107821494683383
$ws.Range("J7").Value = -44.233596

$ws.Range("A8").Value = "~1"
$ws.Range("B8").Value = "~sundur + wind"
$ws.Range("C8").Value = "~1"
$ws.Range("D8").Value = "~1"
$ws.Range("E8").Value = "Phi(~1)p(~sundur + wind)pent(~1)N(~1)"
$ws.Range("F8").Value = 6
$ws.Range("G8").Value = 92.917738
$ws.Range("H8").Value = 5.871139518518518
$ws.Range("I8").Value = 0.01768388706520641
$ws.Range("J8").Value = -41.451341

$ws.Range("A9").Value = "~1"
$ws.Range("B9").Value = "~temp + wind + sundur"
$ws.Range("C9").Value = "~1"
$ws.Range("D9").Value = "~1"
$ws.Range("E9").Value = "Phi(~1)p(~temp + wind + sundur)pent(~1)N(~1)"
$ws.Range("F9").Value = 7
$ws.Range("G9").Value = 93.15929966666667
$ws.Range("H9").Value = 6.112701185185188
$ws.Range("I9").Value = 0.01567195887199801
$ws.Range("J9").Value = -44.516446

$ws.Range("A10").Value = "~1"
$ws.Range("B10").Value = "~1"
$ws.Range("C10").Value = "~time"
$ws.Range("D10").Value = "~1"
$ws.Range("E10").Value = "Phi(~1)p(~1)pent(~time)N(~1)"
$ws.Range("F10").Value = 11
$ws.Range("G10").Value = 100.514298
$ws.Range("H10").Value = 13.46769951851851
$ws.Range("I10").Value = 0.0003962832243656278
$ws.Range("J10").Value = -53.694782

$ws.Range("A11").Value = "~time"
$ws.Range("B11").Value = "~1"
$ws.Range("C11").Value = "~1"
$ws.Range("D11").Value = "~1"
$ws.Range("E11").Value = "Phi(~time)p(~1)pent(~1)N(~1)"
$ws.Range("F11").Value = 11
$ws.Range("G11").Value = 104.21149
$ws.Range("H11").Value = 17.16489151851852
$ws.Range("I11").Value = [double]"6.239799657287439E-05"
$ws.Range("J11").Value = -49.997589

$ws.Range("A12").Value = "~1"
$ws.Range("B12").Value = "~time"
$ws.Range("C12").Value = "~1"
$ws.Range("D12").Value = "~1"
$ws.Range("E12").Value = "Phi(~1)p(~time)pent(~1)N(~1)"
$ws.Range("F12").Value = 12
$ws.Range("G12").Value = 104.7037256315789
$ws.Range("H12").Value = 17.65712715009747
$ws.Range("I12").Value = [double]"4.878463237258471E-05"
$ws.Range("J12").Value = -54.726407

$ws.Range("A13").Value = "~1"
$ws.Range("B13").Value = "~sundur"
$ws.Range("C13").Value = "~time"
$ws.Range("D13").Value = "~1"
$ws.Range("E13").Value = "Phi(~1)p(~sundur)pent(~time)N(~1)"
$ws.Range("F13").Value = 12
$ws.Range("G13").Value = 105.325859631579
$ws.Range("H13").Value = 18.27926115009747
$ws.Range("I13").Value = [double]"3.574278222194166E-05"
$ws.Range("J13").Value = -54.104273

$ws.Range("A14").Value = "~1"
$ws.Range("B14").Value = "~temp"
$ws.Range("C14").Value = "~time"
$ws.Range("D14").Value = "~1"
$ws.Range("E14").Value = "Phi(~1)p(~temp)pent(~time)N(~1)"
$ws.Range("F14").Value = 12
$ws.Range("G14").Value = 105.691611631579
$ws.Range("H14").Value = 18.64501315009747
$ws.Range("I14").Value = [double]"2.976914197662823E-05"
$ws.Range("J14").Value = -53.738521

$ws.Range("A15").Value = "~1"
$ws.Range("B15").Value = "~wind"
$ws.Range("C15").Value = "~time"
$ws.Range("D15").Value = "~1"
$ws.Range("E15").Value = "Phi(~1)p(~wind)pent(~time)N(~1)"
$ws.Range("F15").Value = 12
$ws.Range("G15").Value = 105.7353426315789
$ws.Range("H15").Value = 18.68874415009746
$ws.Range("I15").Value = [double]"2.912528953250178E-05"
$ws.Range("J15").Value = -53.69479

$ws.Range("A16").Value = "~time"
$ws.Range("B16").Value = "~wind"
$ws.Range("C16").Value = "~1"
$ws.Range("D16").Value = "~1"
$ws.Range("E16").Value = "Phi(~time)p(~wind)pent(~1)N(~1)"
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 107.2893546315789
$ws.Range("H16").Value = 20.24275615009746
$ws.Range("I16").Value = [double]"1.339124121898431E-05"
$ws.Range("J16").Value = -52.140778

$ws.Range("A17").Value = "~time"
$ws.Range("B17").Value = "~temp"
$ws.Range("C17").Value = "~1"
$ws.Range("D17").Value = "~1"
$ws.Range("E17").Value = "Phi(~time)p(~temp)pent(~1)N(~1)"
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 108.7790596315789
$ws.Range("H17").Value = 21.73246115009746
$ws.Range("I17").Value = [double]"6.358219293604716E-06"
$ws.Range("J17").Value = -50.651073

$ws.Range("A18").Value = "~time"
$ws.Range("B18").Value = "~sundur"
$ws.Range("C18").Value = "~1"
$ws.Range("D18").Value = "~1"
$ws.Range("E18").Value = "Phi(~time)p(~sundur)pent(~1)N(~1)"
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 109.4321176315789
$ws.Range("H18").Value = 22.38551915009747
$ws.Range("I18").Value = [double]"4.586968520680744E-06"
$ws.Range("J18").Value = -49.998015

$ws.Range("A19").Value = "~1"
$ws.Range("B19").Value = "~sundur + wind"
$ws.Range("C19").Value = "~time"
$ws.Range("D19").Value = "~1"
$ws.Range("E19").Value = "Phi(~1)p(~sundur + wind)pent(~time)N(~1)"
$ws.Range("F19").Value = 13
$ws.Range("G19").Value = 111.0940512222222
$ws.Range("H19").Value = 24.04745274074074
$ws.Range("I19").Value = [double]"1.99821155408167E-06"
$ws.Range("J19").Value = -54.137251

$ws.Range("A20").Value = "~1"
$ws.Range("B20").Value = "~temp + sundur"
$ws.Range("C20").Value = "~time"
$ws.Range("D20").Value = "~1"
$ws.Range("E20").Value = "Phi(~1)p(~temp + sundur)pent(~time)N(~1)"
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 111.1111282222222
$ws.Range("H20").Value = 24.06452974074074
$ws.Range("I20").Value = [double]"1.981222458639957E-06"
$ws.Range("J20").Value = -54.120174

$ws.Range("A21").Value = "~1"
$ws.Range("B21").Value = "~temp + wind"
$ws.Range("C21").Value = "~time"
$ws.Range("D21").Value = "~1"
$ws.Range("E21").Value = "Phi(~1)p(~temp + wind)pent(~time)N(~1)"
$ws.Range("F21").Value = 13
$ws.Range("G21").Value = 111.4904902222222
$ws.Range("H21").Value = 24.44389174074075
$ws.Range("I21").Value = [double]"1.63891273832972E-06"
$ws.Range("J21").Value = -53.740812

$ws.Range("A22").Value = "~time"
$ws.Range("B22").Value = "~temp + wind"
$ws.Range("C22").Value = "~1"
$ws.Range("D22").Value = "~1"
$ws.Range("E22").Value = "Phi(~time)p(~temp + wind)pent(~1)N(~1)"
$ws.Range("F22").Value = 13
$ws.Range("G22").Value = 112.1771802222222
$ws.Range("H22").Value = 25.13058174074075
$ws.Range("I22").Value = [double]"1.162633926627651E-06"
$ws.Range("J22").Value = -53.054122

$ws.Range("A23").Value = "~time"
$ws.Range("B23").Value = "~sundur + wind"
$ws.Range("C23").Value = "~1"
$ws.Range("D23").Value = "~1"
$ws.Range("E23").Value = "Phi(~time)p(~sundur + wind)pent(~1)N(~1)"
$ws.Range("F23").Value = 13
$ws.Range("G23").Value = 112.9398762222222
$ws.Range("H23").Value = 25.89327774074073
$ws.Range("I23").Value = [double]"7.94009429028247E-07"
$ws.Range("J23").Value = -52.291426

$ws.Range("A24").Value = "~time"
$ws.Range("B24").Value = "~temp + sundur"
$ws.Range("C24").Value = "~1"
$ws.Range("D24").Value = "~1"
$ws.Range("E24").Value = "Phi(~time)p(~temp + sundur)pent(~1)N(~1)"
$ws.Range("F24").Value = 13
$ws.Range("G24").Value = 114.4929812222222
$ws.Range("H24").Value = 27.44638274074075
$ws.Range("I24").Value = [double]"3.652356773294304E-07"
$ws.Range("J24").Value = -50.738321

$ws.Range("A25").Value = "~1"
$ws.Range("B25").Value = "~temp + wind + sundur"
$ws.Range("C25").Value = "~time"
$ws.Range("D25").Value = "~1"
$ws.Range("E25").Value = "Phi(~1)p(~temp + wind + sundur)pent(~time)N(~1)"
$ws.Range("F25").Value = 14
$ws.Range("G25").Value = 117.5645023529412
$ws.Range("H25").Value = 30.51790387145969
$ws.Range("I25").Value = [double]"7.863227756384226E-08"
$ws.Range("J25").Value = -54.15046

$ws.Range("A26").Value = "~time"
$ws.Range("B26").Value = "~temp + wind + sundur"
$ws.Range("C26").Value = "~1"
$ws.Range("D26").Value = "~1"
$ws.Range("E26").Value = "Phi(~time)p(~temp + wind + sundur)pent(~1)N(~1)"
$ws.Range("F26").Value = 14
$ws.Range("G26").Value = 118.4601333529412
$ws.Range("H26").Value = 31.41353487145969
$ws.Range("I26").Value = [double]"5.02478003156237E-08"
$ws.Range("J26").Value = -53.254829

$ws.Range("A27").Value = "~time"
$ws.Range("B27").Value = "~1"
$ws.Range("C27").Value = "~time"
$ws.Range("D27").Value = "~1"
$ws.Range("E27").Value = "Phi(~time)p(~1)pent(~time)N(~1)"
$ws.Range("F27").Value = 18
$ws.Range("G27").Value = 148.3504596153846
$ws.Range("H27").Value = 61.30386113390313
$ws.Range("I27").Value = [double]"1.623735087407893E-14"
$ws.Range("J27").Value = -59.274005

$ws.Range("A28").Value = "~time"
$ws.Range("B28").Value = "~wind"
$ws.Range("C28").Value = "~time"
$ws.Range("D28").Value = "~1"
$ws.Range("E28").Value = "Phi(~time)p(~wind)pent(~time)N(~1)"
$ws.Range("F28").Value = 19
$ws.Range("G28").Value = 159.5733643333333
$ws.Range("H28").Value = 72.52676585185185
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = -60.769049

$ws.Range("A29").Value = "~time"
$ws.Range("B29").Value = "~temp"
$ws.Range("C29").Value = "~time"
$ws.Range("D29").Value = "~1"
$ws.Range("E29").Value = "Phi(~time)p(~temp)pent(~time)N(~1)"
$ws.Range("F29").Value = 19
$ws.Range("G29").Value = 160.8222343333333
$ws.Range("H29").Value = 73.77563585185186
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = -59.520179

$ws.Range("A30").Value = "~time"
$ws.Range("B30").Value = "~sundur"
$ws.Range("C30").Value = "~time"
$ws.Range("D30").Value = "~1"
$ws.Range("E30").Value = "Phi(~time)p(~sundur)pent(~time)N(~1)"
$ws.Range("F30").Value = 19
$ws.Range("G30").Value = 160.8692983333333
$ws.Range("H30").Value = 73.82269985185187
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = -59.473115

$ws.Range("A31").Value = "~time"
$ws.Range("B31").Value = "~time"
$ws.Range("C31").Value = "~1"
$ws.Range("D31").Value = "~1"
$ws.Range("E31").Value = "Phi(~time)p(~time)pent(~1)N(~1)"
$ws.Range("F31").Value = 19
$ws.Range("G31").Value = 160.9105593333333
$ws.Range("H31").Value = 73.86396085185186
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = -59.431854

$ws.Range("A32").Value = "~1"
$ws.Range("B32").Value = "~time"
$ws.Range("C32").Value = "~time"
$ws.Range("D32").Value = "~1"
$ws.Range("E32").Value = "Phi(~1)p(~time)pent(~time)N(~1)"
$ws.Range("F32").Value = 19
$ws.Range("G32").Value = 161.2669543333333
$ws.Range("H32").Value = 74.22035585185185
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = -59.075459

$ws.Range("A33").Value = "~time"
$ws.Range("B33").Value = "~temp + wind"
$ws.Range("C33").Value = "~time"
$ws.Range("D33").Value = "~1"
$ws.Range("E33").Value = "Phi(~time)p(~temp + wind)pent(~time)N(~1)"
$ws.Range("F33").Value = 20
$ws.Range("G33").Value = 172.5183543636364
$ws.Range("H33").Value = 85.47175588215489
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = -62.854361

$ws.Range("A34").Value = "~time"
$ws.Range("B34").Value = "~sundur + wind"
$ws.Range("C34").Value = "~time"
$ws.Range("D34").Value = "~1"
$ws.Range("E34").Value = "Phi(~time)p(~sundur + wind)pent(~time)N(~1)"
$ws.Range("F34").Value = 20
$ws.Range("G34").Value = 174.5819133636364
$ws.Range("H34").Value = 87.53531488215489
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = -60.790803

$ws.Range("A35").Value = "~time"
$ws.Range("B35").Value = "~temp + sundur"
$ws.Range("C35").Value = "~time"
$ws.Range("D35").Value = "~1"
$ws.Range("E35").Value = "Phi(~time)p(~temp + sundur)pent(~time)N(~1)"
$ws.Range("F35").Value = 20
$ws.Range("G35").Value = 175.3834563636364
$ws.Range("H35").Value = 88.33685788215487
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = -59.98926

$ws.Range("A36").Value = "~time"
$ws.Range("B36").Value = "~temp + wind + sundur"
$ws.Range("C36").Value = "~time"
$ws.Range("D36").Value = "~1"
$ws.Range("E36").Value = "Phi(~time)p(~temp + wind + sundur)pent(~time)N(~1)"
$ws.Range("F36").Value = 21
$ws.Range("G36").Value = 190.55124
$ws.Range("H36").Value = 103.5046415185185
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = -62.85784

$ws.Range("A37").Value = "~time"
$ws.Range("B37").Value = "~time"
$ws.Range("C37").Value = "~time"
$ws.Range("D37").Value = "~1"
$ws.Range("E37").Value = "Phi(~time)p(~time)pent(~time)N(~1)"
$ws.Range("F37").Value = 26
$ws.Range("G37").Value = 386.937527
$ws.Range("H37").Value = 299.8909285185185
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = -64.87155300000001

